# PHOENIX-6061: changes in Grievances
# Adds two new login-test-data rows (CSCUser / PublicHealthJA) to the
# registeredUserDetails sheet, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: CSCUser / CSCUSER / kurnool_eGov@123 (hyperlink) / FALSE() formula
$ws.Range("A31").Value = "CSCUser"
$ws.Range("B31").Value = "CSCUSER"
$ws.Range("C31").Value = "kurnool_eGov@123"
$ws.Hyperlinks.Add($ws.Range("C31"), "mailto:kurnool_eGov@123") | Out-Null
$ws.Range("D31").Formula = "=FALSE()"

# Row 32: PublicHealthJA / 0944182 (kept as text, leading zero) / kurnool_eGov@123 / TRUE
$ws.Range("A32").Value = "PublicHealthJA"
$ws.Range("B32").Value = "'0944182"
$ws.Range("C32").Value = "kurnool_eGov@123"
$ws.Range("D32").Value = $true

# Restore the plain "Explanatory Text" cell style (General format) on the new
# rows - typing the leading apostrophe above nudges Excel into the
# "quote prefix" style, so reapply the base style used by the rest of the data.
$ws.Range("A31:D32").Style = "Explanatory Text"

# Move the view/selection the way the author left it.
$ws.Range("A35").Select()
